$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("FCFS")

# Update the FCFS burst times (column C) -- the Wait/Turnaround/Response
# formulas in D:F recalc automatically from these.
$ws1.Range("C2").Value = 8
$ws1.Range("C4").Value = 9
$ws1.Range("C5").Value = 5

# Duplicate the FCFS sheet (keeps formulas, formatting, merged cells) and
# place the copy right after FCFS, then rename it to SJF.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1) | Out-Null
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "SJF"

# SJF hasn't had its burst times filled in yet -- reset column C to 0.
$ws2.Range("C2:C6").Value = 0

# Restore each sheet's own selection; selecting SJF last leaves it as the
# active sheet/tab, matching the saved workbook state.
$ws1.Range("C7").Select() | Out-Null
$ws2.Range("C9").Select() | Out-Null
